$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(86).Insert()

$ws.Range("C86").Value = "Scalpel Accuracy:"
$ws.Range("D86").Value = 492.86
